$p = $ppt.ActivePresentation

# --- Change 1: Table 6 on slide 11 - add @start_time and @end_time paragraphs ---
$s11 = $p.Slides.Item(11)
$shp6 = $s11.Shapes.Item(5)
$cell = $shp6.Table.Cell(2, 1)
$tr = $cell.Shape.TextFrame.TextRange
$tr.Text = "@id`r@budget(numeric)`r@start_time`r@end_time"

# --- Change 2: Straight Arrow Connector 19 on slide 9 - reposition + flip vertical ---
$s9 = $p.Slides.Item(9)
$conn = $s9.Shapes.Item(10)
$conn.Left = 465.5294488
$conn.Top = 148.0
$conn.Width = 59.80385827
$conn.Height = 34.30755906
$conn.Flip(1)
